$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.000000079193403877860862803177
$ws.Range("E2").Value = 0.000000079193403877860862803177

$ws.Range("D3").Value = 0.999999967968157599251810552232
$ws.Range("E3").Value = 0.999999967968157599251810552232

$ws.Range("D4").Value = 0.995580763019970671479086377076
$ws.Range("E4").Value = 0.004419236980029329388275360913

$ws.Range("D5").Value = 0.999999999999971578290569595993
$ws.Range("E5").Value = 0.000000000000028421709430404007

$ws.Range("D6").Value = 0.003450676946001859056456506636
$ws.Range("E6").Value = 0.996549323053998126198393947561

$ws.Range("D7").Value = 0.999999987400975842710693086701
$ws.Range("E7").Value = 0.000000012599024157289310222021

$ws.Range("D8").Value = 0.002634377503377586821758926661
$ws.Range("E8").Value = 0.997365622496622417081368894287
$ws.Range("F8").Value = 4.124179363250732421875000000000
